$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "FIXED $ OR $/UNIT/PERIOD" values between the Wet Storage (row 2)
# and Dry Storage (row 6) records.
$ws.Range("F2").Value = "Wet Storage"
$ws.Range("F6").Value = "Dry Storage"

# The CAPACITY / PRICING TYPE columns (C:D) no longer carry the extra
# applied-font style - clear formatting back to the default style.
$ws.Range("C2:D6").ClearFormats()

# Give column E (blank spacer column) an explicit custom width.
$ws.Columns("E:E").ColumnWidth = 9.7

# Move/update the active selection.
$ws.Range("F15").Select() | Out-Null
